# Adds two new per-device result sheets ("iPhone-6 Avner" and "iPhone-6 null"),
# mirroring the layout/content of the existing "iPhone-6 Brian"/"iPhone-6 Raj"
# sheets, including the FAIL-cell highlighting and the external screenshot
# hyperlinks on each FAIL cell.

$wb = $excel.ActiveWorkbook

$paramsHeader = "Test Parameters"
$sasha        = "testSignIn, shshshs@aaa.com, dsdsdsd, Hi, Sasha!"
$saTeam       = "testSignIn, SATeam@perfectomobile.com, SATeam123, Hi, SA!"
$avner        = "testSignIn, avnerg@perfectomobile.com, a1001a, Hi, Avner!"
$fail         = "FAIL"

# ---------------------------------------------------------------------
# Sheet: "iPhone-6 Avner"  (two result columns: B = build 3.05, C = 3.06)
# ---------------------------------------------------------------------
$lastIdx = $wb.Worksheets.Count
$lastWs = $wb.Worksheets.Item($lastIdx)
$wsAvner = $wb.Worksheets.Add($null, $lastWs)
$wsAvner.Name = "iPhone-6 Avner"

$wsAvner.Range("A1").Value = $paramsHeader
$wsAvner.Range("B1").Value = "Avner 3.05"
$wsAvner.Range("C1").Value = "Avner 3.06"

$wsAvner.Range("A2").Value = $sasha
$wsAvner.Range("A3").Value = $saTeam
$wsAvner.Range("A4").Value = $avner

$avnerLinks = @{
    "B2" = "C:/Users/rajp/git/Beton/Beton/test-output/screenshots/2015-08-04-10-22-18-EDT.png.png"
    "B3" = "C:/Users/rajp/git/Beton/Beton/test-output/screenshots/2015-08-04-10-24-26-EDT.png.png"
    "B4" = "C:/Users/rajp/git/Beton/Beton/test-output/screenshots/2015-08-04-10-26-33-EDT.png.png"
    "C2" = "C:/Users/rajp/git/Beton/Beton/test-output/screenshots/2015-08-04-10-28-41-EDT.png.png"
    "C3" = "C:/Users/rajp/git/Beton/Beton/test-output/screenshots/2015-08-04-10-30-49-EDT.png.png"
    "C4" = "C:/Users/rajp/git/Beton/Beton/test-output/screenshots/2015-08-04-10-32-57-EDT.png.png"
}

foreach ($addr in @("B2", "B3", "B4", "C2", "C3", "C4")) {
    $cell = $wsAvner.Range($addr)
    $cell.Value = $fail
    $cell.Interior.Color = 255
    $wsAvner.Hyperlinks.Add($cell, $avnerLinks[$addr])
    # Keep the plain FAIL-fill look (no Excel auto "Hyperlink" blue/underline font)
    $cell.Font.Underline = $false
    $cell.Font.ColorIndex = 1
}

# ---------------------------------------------------------------------
# Sheet: "iPhone-6 null"  (single result column: B = build 3.07)
# ---------------------------------------------------------------------
$lastIdx2 = $wb.Worksheets.Count
$lastWs2 = $wb.Worksheets.Item($lastIdx2)
$wsNull = $wb.Worksheets.Add($null, $lastWs2)
$wsNull.Name = "iPhone-6 null"

$wsNull.Range("A1").Value = $paramsHeader
$wsNull.Range("B1").Value = "Avner 3.07"

$wsNull.Range("A2").Value = $sasha
$wsNull.Range("A3").Value = $saTeam
$wsNull.Range("A4").Value = $avner

$nullLinks = @{
    "B2" = "C:/Users/rajp/git/Beton/Beton/test-output/screenshots/2015-08-04-10-35-05-EDT.png.png"
    "B3" = "C:/Users/rajp/git/Beton/Beton/test-output/screenshots/2015-08-04-10-37-13-EDT.png.png"
    "B4" = "C:/Users/rajp/git/Beton/Beton/test-output/screenshots/2015-08-04-10-39-21-EDT.png.png"
}

foreach ($addr in @("B2", "B3", "B4")) {
    $cell = $wsNull.Range($addr)
    $cell.Value = $fail
    $cell.Interior.Color = 255
    $wsNull.Hyperlinks.Add($cell, $nullLinks[$addr])
    $cell.Font.Underline = $false
    $cell.Font.ColorIndex = 1
}

# Restore original active sheet/tab selection (unchanged by this edit).
$wb.Worksheets.Item(1).Activate()
